# Add Week 12 betting lines to the "Sheet1" worksheet (the active sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# week, game, total_line, spread_line
$week12 = @(
    @("BUF_HOU", 46.5, -3.5),
    @("PIT_CHI", 44.5, 2.5),
    @("NYJ_BAL", 44.5, 11.5),
    @("NYG_DET", 45.5, 8.5),
    @("NE_CIN", 47.5, 4.5),
    @("SEA_TEN", 44.5, -1.5),
    @("MIN_GB", 45.5, 3.5),
    @("IND_KC", 45.5, 9.5),
    @("JAX_ARI", 47.5, 3),
    @("CLE_LV", 41.5, 3.5),
    @("ATL_NO", 43.5, -2.5),
    @("PHI_DAL", 46.5, -4.5),
    @("TB_LA", 47.5, 1.5),
    @("CAR_SF", 46.5, 7)
)

$startRow = 150
$r = $startRow
foreach ($game in $week12) {
    $ws.Cells.Item($r, 1).Value = 12
    $ws.Cells.Item($r, 2).Value = $game[0]
    $ws.Cells.Item($r, 3).Value = $game[1]
    $ws.Cells.Item($r, 4).Value = $game[2]
    $r = $r + 1
}

# Update the window scroll position / selection to match where the author
# was working after adding the new rows.
[void]$excel.Goto($ws.Range("A142"), $false)
[void]$ws.Range("C149").Select()
